# Applies the "change position of area_pv" edit:
#  - ELECTRICITY sheet loses its area_pv / efficiency_pv / CAPEX_PV_USD2015kW
#    columns (they move out into a brand-new ELECTRICITY_PV sheet).
#  - A new ELECTRICITY_PV sheet is added after ELECTRICITY with
#    Description/code/component/area/reference columns describing the two
#    PV assemblies.

$wb = $excel.ActiveWorkbook

$electricity = $wb.Worksheets.Item("ELECTRICITY")
$heating = $wb.Worksheets.Item("HEATING")

# --- 1. Remove the PV-specific columns (area_pv, efficiency_pv,
#        CAPEX_PV_USD2015kW) from ELECTRICITY; LT_yr / O&M_% shift left
#        into columns G:H.
$electricity.Columns("G:I").Delete()

# --- 2. Add the new ELECTRICITY_PV sheet right after ELECTRICITY.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$pv = $wb.Worksheets.Add($null, $lastSheet)
$pv.Name = "ELECTRICITY_PV"

# --- 3. Copy cell formatting (styles) from existing cells that already
#        carry the exact formats we need, so new cells match the
#        workbook's established look (borders/fonts/number formats).

# Header row (Description/code/component/area/reference) - same style as
# the ELECTRICITY header row.
$electricity.Range("A1:E1").Copy()
$pv.Range("A1").PasteSpecial(-4122)

# Row 2 ("none" baseline assembly).
$electricity.Range("A2:C2").Copy()
$pv.Range("A2").PasteSpecial(-4122)
$heating.Range("E2").Copy()
$pv.Range("D2").PasteSpecial(-4122)
$heating.Range("J2").Copy()
$pv.Range("E2").PasteSpecial(-4122)

# Row 3 (PV1 assembly).
$electricity.Range("A3:C3").Copy()
$pv.Range("A3").PasteSpecial(-4122)
$heating.Range("F3").Copy()
$pv.Range("D3").PasteSpecial(-4122)
$heating.Range("J2").Copy()
$pv.Range("E3").PasteSpecial(-4122)

# --- 4. Fill in the values.
$pv.Range("A1").Value = "Description"
$pv.Range("B1").Value = "code"
$pv.Range("C1").Value = "component"
$pv.Range("D1").Value = "area"
$pv.Range("E1").Value = "reference"

$pv.Range("A2").Value = "none"
$pv.Range("B2").Value = "SUPPLY_ELECTRICITY_PV_AS0"
$pv.Range("C2").Value = "NONE"
$pv.Range("D2").Value = 0
$pv.Range("E2").Value = "educated guess"

$pv.Range("A3").Value = "20 m2 Fotovoltaic installation Moncrytalline"
$pv.Range("B3").Value = "SUPPLY_ELECTRICITY_PV_AS1"
$pv.Range("C3").Value = "PV1"
$pv.Range("D3").Value = 20
$pv.Range("E3").Value = "educated guess"

Write-Host "Edit applied"
